$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1 = Wins, AE1 = Losses, AF1 = Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header cells (e.g. A1) so the new
# header cells share the same style as the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-45: team record Wins=84, Losses=78, Ties=0
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 84  # AD
    $ws.Cells.Item($r, 31).Value = 78  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
